$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.327.63"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "'1.814.52"
$ws.Range("E3").Value = "  +3.67%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'326.97"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").Value = "'0.4371"
$ws.Range("E7").Value = "  +2.43%  "

$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("D9").Value = "'44.68"
$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("D10").Value = "'0.07675"
$ws.Range("E10").Value = "  +2.69%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "'22.00"
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").Value = "'6.316"
$ws.Range("E14").Value = "  +3.16%  "

$ws.Range("D15").Value = "'7.520"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").Value = "'1.820.39"
$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("D17").Value = "'95.46"
$ws.Range("E17").Value = "  +8.66%  "

$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").Value = "'0.06512"
$ws.Range("E19").Value = "  +4.40%  "

$ws.Range("D20").Value = "'0.9997"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("E21").Value = "  +2.66%  "

$ws.Range("D22").Value = "'6.247"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").Value = "'28.337.20"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").Value = "'11.58"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").Value = "'2.093"
$ws.Range("E25").Value = "  -9.98%  "

$ws.Range("D26").Value = "'162.16"
$ws.Range("E26").Value = "  +6.88%  "

$ws.Range("D27").Value = "'20.72"
$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("D28").Value = "'2.031.78"
$ws.Range("E28").Value = "  +4.67%  "

$ws.Range("D29").Value = "'2.276"
$ws.Range("E29").Value = "  -4.02%  "

$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").Value = "'1.207"
$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("D32").Value = "'5.974"
$ws.Range("E32").Value = "  +4.43%  "

$ws.Range("D33").Value = "'0.09186"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "'3.513"
$ws.Range("E34").Value = "  -4.62%  "

$ws.Range("D35").Value = "'13.00"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("D36").Value = "'0.02349"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2173"
$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.197"
$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("D39").Value = "'0.6583"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").Value = "'0.06208"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("D41").Value = "'1.195"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "'8.117"
$ws.Range("E42").Value = "  +2.14%  "

$ws.Range("D43").Value = "'1.429"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").Value = "'13.86"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "'0.6113"
$ws.Range("E46").Value = "  +3.45%  "

$ws.Range("D47").Value = "'3.748"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("D48").Value = "'125.94"
$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").Value = "'2.016"
$ws.Range("E49").Value = "  +2.66%  "

$ws.Range("E50").Value = "  +2.79%  "

$ws.Range("D51").Value = "'0.07003"
$ws.Range("E51").Value = "  +1.76%  "
